# Update "Prix Spot" sheet: add a new day column AE (14-jul) with hourly prices.
$wb = $excel.ActiveWorkbook

$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting (bold, centered, bordered) from the previous header cell
# so the new header cell AE1 matches the style of the existing day headers.
$wsPrix.Range("AD1").Copy() | Out-Null
$wsPrix.Range("AE1").PasteSpecial(-4122) | Out-Null
$wsPrix.Range("AE1").Value = "14-jul"

$prixValues = @{
    2  = 85.26000000000001
    3  = 71.56
    4  = 49.28
    5  = 51.95
    6  = 41.39
    7  = 35.78
    8  = 40.19
    9  = 31.17
    10 = 28.27
    11 = 38.53
    12 = 37.05
    13 = 30.34
    14 = 34.06
    15 = 38.8
    16 = 31.49
    17 = 18.21
    18 = 27.09
    19 = 47.5
    20 = 54.14
    21 = 77.76000000000001
    22 = 86.23
    23 = 89.29000000000001
    24 = 117.84
    25 = 113.83
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 31).Value = $prixValues[$row]
}

# Helper: write a literal (non-date) text value into a cell. Plain
# assignment of a date-shaped string (e.g. "2025-07-12") would otherwise be
# auto-converted to a date serial by Excel's normal text-entry parsing, so
# the cell is temporarily forced to Text format for the write and restored
# to the default "Normal" style afterwards (matches the plain string cells
# already used for the other dates in these columns).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Update "Gaz" sheet: append two new daily rows.
$wsGaz = $wb.Worksheets.Item("Gaz")
Set-TextValue $wsGaz.Range("A28") "2025-07-12"
$wsGaz.Range("B28").Value = 34.8
Set-TextValue $wsGaz.Range("A29") "2025-07-13"
$wsGaz.Range("B29").Value = 34.8

# Update "CO2" sheet: append two new daily rows.
$wsCo2 = $wb.Worksheets.Item("CO2")
Set-TextValue $wsCo2.Range("A28") "2025-07-12"
$wsCo2.Range("B28").Value = 69.8
Set-TextValue $wsCo2.Range("A29") "2025-07-13"
$wsCo2.Range("B29").Value = 69.8
